# Update countries & provincias Spain
# Refresh the COVID snapshot figures and re-rank a couple of
# countries whose "Casos totales" crossed each other.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- timestamp header (A1) -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 30 de Julio de 2020 a las 21:11"

# --- Estados Unidos (row 4) -------------------------------------------------
$ws.Range("B4").Value = 4607701
$ws.Range("C4").Value = 39664
$ws.Range("D4").Value = 2255763
$ws.Range("E4").Value = 2197318
$ws.Range("G4").Value = 780
$ws.Range("H4").Value = 154620

# --- India (row 6) -----------------------------------------------------------
$ws.Range("B6").Value = 1639350
$ws.Range("C6").Value = 54966
$ws.Range("D6").Value = 1059093
$ws.Range("E6").Value = 544471

# --- Reino Unido (row 13) ----------------------------------------------------
$ws.Range("G13").Value = 38
$ws.Range("H13").Value = 45999

# --- Alemania (row 21) -------------------------------------------------------
$ws.Range("B21").Value = 209601
$ws.Range("C21").Value = 790
$ws.Range("E21").Value = 8380
$ws.Range("G21").Value = 9
$ws.Range("H21").Value = 9221

# --- Francia (row 22) --------------------------------------------------------
$ws.Range("B22").Value = 186573
$ws.Range("C22").Value = 1377
$ws.Range("E22").Value = 75008
$ws.Range("G22").Value = 16
$ws.Range("H22").Value = 30254

# --- China / Ecuador swap places (rows 31-32) --------------------------------
# Ecuador's updated total (84370) now exceeds China's static total (84165),
# so Ecuador moves up to row 31 and China drops to row 32.
$ws.Range("A31").Value = "Ecuador"
$ws.Range("B31").Value = 84370
$ws.Range("C31").Value = 1177
$ws.Range("D31").Value = 35824
$ws.Range("E31").Value = 42889
$ws.Range("G31").Value = 34
$ws.Range("H31").Value = 5657

$ws.Range("A32").Value = "China"
$ws.Range("B32").Value = 84165
$ws.Range("C32").Value = 105
$ws.Range("D32").Value = 78957
$ws.Range("E32").Value = 574
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 4634

# --- Guatemala (row 48) ------------------------------------------------------
$ws.Range("B48").Value = 48826
$ws.Range("C48").Value = 1221
$ws.Range("D48").Value = 35629
$ws.Range("E48").Value = 11330
$ws.Range("G48").Value = 32
$ws.Range("H48").Value = 1867

# --- Estonia / Namibia swap places (rows 127-128) ----------------------------
# Namibia's updated total (2052) now exceeds Estonia's static total (2051),
# so Namibia moves up to row 127 and Estonia drops to row 128.
$ws.Range("A127").Value = "Namibia"
$ws.Range("B127").Value = 2052
$ws.Range("C127").Value = 66
$ws.Range("D127").Value = 164
$ws.Range("E127").Value = 1878
$ws.Range("G127").Value = 1
$ws.Range("H127").Value = 10

$ws.Range("A128").Value = "Estonia"
$ws.Range("B128").Value = 2051
$ws.Range("C128").Value = 9
$ws.Range("D128").Value = 1926
$ws.Range("E128").Value = 56
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 69

# --- Republica de Chipre (row 147) ------------------------------------------
$ws.Range("B147").Value = 1090
$ws.Range("C147").Value = 10
$ws.Range("E147").Value = 219
